$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 72 (this shifts the existing rows 72..164
# down to 73..165, matching the diff's row-shift pattern).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new data record.
$ws.Range("A72").Value = 3
$ws.Range("B72").Value = "Femacal de La Calera"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 44413
$ws.Range("E72").Value = 5
$ws.Range("F72").Value = 100114013
$ws.Range("G72").Value = "Zanahoria"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 250
$ws.Range("K72").Value = 5000
$ws.Range("L72").Value = 5000
$ws.Range("M72").Value = 5000
$ws.Range("N72").Value = "$/saco 20 kilos"
$ws.Range("O72").Value = "Provincia de Quillota"
$ws.Range("P72").Value = 250
$ws.Range("Q72").Value = 20
$ws.Range("R72").Value = "Hortaliza"
